# Update countries & provincias Spain
# Refresh COVID-19 country statistics table and reorder a handful of
# countries whose case totals changed their relative ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (footer timestamp)
$ws.Range("A1").Value = 'Datos actualizados a 11 de Abril de 2020 a las 19:22'

# Row 4
$ws.Range("B4").Value = 520999
$ws.Range("C4").Value = 18123
$ws.Range("D4").Value = 28507
$ws.Range("E4").Value = 472357
$ws.Range("G4").Value = 1388
$ws.Range("H4").Value = 20135

# Row 12
$ws.Range("B12").Value = 52167
$ws.Range("C12").Value = 5138
$ws.Range("D12").Value = 2965
$ws.Range("E12").Value = 48101
$ws.Range("F12").Value = 1626
$ws.Range("G12").Value = 95
$ws.Range("H12").Value = 1101

# Row 14
$ws.Range("B14").Value = 25071
$ws.Range("C14").Value = 520
$ws.Range("E14").Value = 12935
$ws.Range("G14").Value = 34
$ws.Range("H14").Value = 1036

# Row 16
$ws.Range("B16").Value = 22575
$ws.Range("C16").Value = 427
$ws.Range("E16").Value = 15993

# Row 19
$ws.Range("B19").Value = 13798
$ws.Range("C19").Value = 238
$ws.Range("E19").Value = 6857

# Row 24
$ws.Range("A24").Value = 'India'
$ws.Range("B24").Value = 8339
$ws.Range("C24").Value = 739
$ws.Range("D24").Value = 774
$ws.Range("E24").Value = 7316
$ws.Range("F24").Value = 0
$ws.Range("H24").Value = 249

# Row 25
$ws.Range("A25").Value = 'Irlanda'
$ws.Range("B25").Value = 8089
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 7777
$ws.Range("F25").Value = 194
$ws.Range("H25").Value = 287

# Row 33
$ws.Range("E33").Value = 4942
$ws.Range("G33").Value = 20
$ws.Range("H33").Value = 290

# Row 55
$ws.Range("B55").Value = 2028
$ws.Range("C55").Value = 25
$ws.Range("E55").Value = 1593
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 25

# Row 81
$ws.Range("B81").Value = 685
$ws.Range("C81").Value = 14
$ws.Range("E81").Value = 614
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 28

# Row 86
$ws.Range("D86").Value = 61
$ws.Range("E86").Value = 545
$ws.Range("F86").Value = 8

# Row 101
$ws.Range("B101").Value = 381
$ws.Range("C101").Value = 9
$ws.Range("D101").Value = 177
$ws.Range("E101").Value = 197

# Row 133
$ws.Range("B133").Value = 92
$ws.Range("C133").Value = 2
$ws.Range("E133").Value = 86

# Row 134
$ws.Range("A134").Value = 'Aruba'
$ws.Range("B134").Value = 92
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 29
$ws.Range("E134").Value = 63
$ws.Range("H134").Value = 0

# Row 135
$ws.Range("A135").Value = 'Mali'
$ws.Range("B135").Value = 87
$ws.Range("D135").Value = 22
$ws.Range("E135").Value = 58
$ws.Range("H135").Value = 7

# Row 146
$ws.Range("A146").Value = 'Liberia'
$ws.Range("C146").Value = 11
$ws.Range("D146").Value = 3
$ws.Range("E146").Value = 40
$ws.Range("F146").Value = 0
$ws.Range("H146").Value = 5

# Row 147
$ws.Range("A147").Value = 'Bermudas'
$ws.Range("B147").Value = 48
$ws.Range("D147").Value = 25
$ws.Range("E147").Value = 19
$ws.Range("F147").Value = 2
$ws.Range("H147").Value = 4

# Row 148
$ws.Range("A148").Value = 'Gabon'
$ws.Range("B148").Value = 46
$ws.Range("C148").Value = 2
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 44
$ws.Range("F148").Value = 0
$ws.Range("H148").Value = 1

# Row 149
$ws.Range("A149").Value = 'Islas Caimanes'
$ws.Range("B149").Value = 45
$ws.Range("D149").Value = 6
$ws.Range("E149").Value = 38

# Row 150
$ws.Range("A150").Value = 'Macao'
$ws.Range("B150").Value = 45
$ws.Range("D150").Value = 10
$ws.Range("E150").Value = 35
$ws.Range("H150").Value = 0

# Row 151
$ws.Range("A151").Value = 'Bahamas'
$ws.Range("B151").Value = 42
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 5
$ws.Range("E151").Value = 29
$ws.Range("F151").Value = 1
$ws.Range("H151").Value = 8

# Row 152
$ws.Range("A152").Value = 'Guyana'
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 8
$ws.Range("E152").Value = 26
$ws.Range("F152").Value = 3
$ws.Range("H152").Value = 6

# Row 153
$ws.Range("A153").Value = 'Zambia'
$ws.Range("B153").Value = 40
$ws.Range("D153").Value = 28
$ws.Range("E153").Value = 10
$ws.Range("F153").Value = 1

# Row 154
$ws.Range("A154").Value = 'Puerto Rico'
$ws.Range("B154").Value = 39
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 1
$ws.Range("E154").Value = 36
$ws.Range("H154").Value = 2

# Row 155
$ws.Range("A155").Value = 'Guinea-Bisau'
$ws.Range("B155").Value = 38
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 38
$ws.Range("H155").Value = 0

# Row 163
$ws.Range("D163").Value = 5
$ws.Range("E163").Value = 15
$ws.Range("G163").Value = 3
$ws.Range("H163").Value = 5

# Row 197
$ws.Range("A197").Value = 'Islas Turcas y Caicos'
$ws.Range("C197").Value = 0

# Row 198
$ws.Range("A198").Value = 'Nicaragua'
$ws.Range("C198").Value = 1

# Row 203
$ws.Range("A203").Value = 'Burundi'
$ws.Range("C203").Value = 2
$ws.Range("D203").Value = 0
$ws.Range("E203").Value = 5

# Row 204
$ws.Range("A204").Value = 'Islas Malvinas'
$ws.Range("D204").Value = 1
$ws.Range("E204").Value = 4

# Row 205
$ws.Range("A205").Value = 'Butan'
$ws.Range("B205").Value = 5
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 3

# Row 206
$ws.Range("A206").Value = 'Santo Tome y Principe'

# Row 207
$ws.Range("A207").Value = 'Sahara Occidental'

# Row 208
$ws.Range("A208").Value = 'Sudan del Sur'
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 4
$ws.Range("H208").Value = 0

# Row 209
$ws.Range("A209").Value = 'Gambia'
$ws.Range("B209").Value = 4
$ws.Range("D209").Value = 2
$ws.Range("E209").Value = 1
$ws.Range("H209").Value = 1
